$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: WALMART / Tenth Order, amount + invoice no changed
$ws.Range("C2").Value = 30000
$ws.Range("D2").Value = 1109
$ws.Range("E2").Value = "WALMART"
$ws.Range("F2").Value = "Tenth Order"

# Row 3: HOMEDEP (with wrap-text style copied from old E4/E5), amount + invoice no changed
$ws.Range("C3").Value = 35000
$ws.Range("D3").Value = 1106
$ws.Range("E4").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "HOMEDEP"
$ws.Range("F3").Value = "Tenth Order"

# Row 4: clear the amount/invoice/description cells, keep dates+customer blank but formatted
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()

# Row 5: delete entirely, shifting rows up
$ws.Rows("5").Delete()
